$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 663
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 550.6667
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 550.6667
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -900.6667

$ws.Range("H33").Value = 694.15
$ws.Range("I33").Value = 210.64706
$ws.Range("K33").Value = 210.64706
$ws.Range("M33").Value = 18.35293999999999

$ws.Range("H129").Value = 527344.5
$ws.Range("J129").Value = 667904.6
$ws.Range("L129").Value = 2003713.8
$ws.Range("N129").Value = -2013713.8

$ws.Range("H132").Value = 37040564
$ws.Range("I132").Value = 43482144
$ws.Range("J132").Value = 1476.5
$ws.Range("K132").Value = 130446432
$ws.Range("L132").Value = 4429.5
$ws.Range("M132").Value = -130443902
$ws.Range("N132").Value = -9489.5

$ws.Range("H137").Value = 169296.58
$ws.Range("I137").Value = 269089.06
$ws.Range("J137").Value = 2975.7778
$ws.Range("K137").Value = 807267.1799999999
$ws.Range("L137").Value = 8927.3334
$ws.Range("M137").Value = -804717.1799999999
$ws.Range("N137").Value = -14027.3334

$ws.Range("H138").Value = 3143.1628
$ws.Range("I138").Value = 2209.3333
$ws.Range("J138").Value = 3815.52
$ws.Range("K138").Value = 6627.999899999999
$ws.Range("L138").Value = 11446.56
$ws.Range("M138").Value = -1487.999899999999
$ws.Range("N138").Value = -21726.56

$ws.Range("H141").Value = 2485.276
$ws.Range("I141").Value = 2137.423
$ws.Range("K141").Value = 6412.268999999999
$ws.Range("M141").Value = -1232.268999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9504.781999999999
$ws.Range("I32").Value = 7336.902
$ws.Range("J32").Value = 26034.875
$ws.Range("K32").Value = 7336.902
$ws.Range("L32").Value = 26034.875
$ws.Range("M32").Value = -7049.902
$ws.Range("N32").Value = -26608.875

$ws.Range("H61").Value = 10104417
$ws.Range("I61").Value = 12824041
$ws.Range("J61").Value = 2957
$ws.Range("K61").Value = 12824041
$ws.Range("L61").Value = 2957
$ws.Range("M61").Value = -12823829
$ws.Range("N61").Value = -3381

$ws.Range("H74").Value = 32259734
$ws.Range("I74").Value = 55556300
$ws.Range("J74").Value = 2953.077
$ws.Range("K74").Value = 55556300
$ws.Range("L74").Value = 2953.077
$ws.Range("M74").Value = -55555426
$ws.Range("N74").Value = -4701.077

$ws.Range("H77").Value = 32259734
$ws.Range("I77").Value = 55556300
$ws.Range("J77").Value = 2953.077
$ws.Range("K77").Value = 277781500
$ws.Range("L77").Value = 14765.385
$ws.Range("M77").Value = -277777132
$ws.Range("N77").Value = -23501.385

$ws.Range("H132").Value = 9270696
$ws.Range("I132").Value = 11906674
$ws.Range("J132").Value = 44770.832
$ws.Range("K132").Value = 35720022
$ws.Range("L132").Value = 134312.496
$ws.Range("M132").Value = -35717492
$ws.Range("N132").Value = -139372.496

$ws.Range("H136").Value = 10104417
$ws.Range("I136").Value = 12824041
$ws.Range("J136").Value = 2957
$ws.Range("K136").Value = 38472123
$ws.Range("L136").Value = 8871
$ws.Range("M136").Value = -38469573
$ws.Range("N136").Value = -13971

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4991.457
$ws.Range("J134").Value = 4209.6665
$ws.Range("L134").Value = 12628.9995
$ws.Range("N134").Value = -17698.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 153.14285
$ws.Range("I22").Value = 148
$ws.Range("J22").Value = 175
$ws.Range("K22").Value = 148
$ws.Range("L22").Value = 175
$ws.Range("M22").Value = 202
$ws.Range("N22").Value = -875

$ws.Range("H58").Value = 14951
$ws.Range("I58").Value = 1688.2858
$ws.Range("J58").Value = 23023.957
$ws.Range("K58").Value = 1688.2858
$ws.Range("L58").Value = 23023.957
$ws.Range("M58").Value = -1485.2858
$ws.Range("N58").Value = -23429.957

$ws.Range("H99").Value = 3335.6775
$ws.Range("I99").Value = 2522.24
$ws.Range("J99").Value = 6725
$ws.Range("K99").Value = 2522.24
$ws.Range("L99").Value = 6725
$ws.Range("M99").Value = -1024.24
$ws.Range("N99").Value = -9721

$ws.Range("H126").Value = 3335.6775
$ws.Range("I126").Value = 2522.24
$ws.Range("J126").Value = 6725
$ws.Range("K126").Value = 7566.719999999999
$ws.Range("L126").Value = 20175
$ws.Range("M126").Value = -5096.719999999999
$ws.Range("N126").Value = -25115

$ws.Range("H134").Value = 41667536
$ws.Range("I134").Value = 47619870
$ws.Range("J134").Value = 1166.6666
$ws.Range("K134").Value = 142859610
$ws.Range("L134").Value = 3499.9998
$ws.Range("M134").Value = -142857075
$ws.Range("N134").Value = -8569.9998

$ws.Range("H136").Value = 14951
$ws.Range("I136").Value = 1688.2858
$ws.Range("J136").Value = 23023.957
$ws.Range("K136").Value = 5064.857400000001
$ws.Range("L136").Value = 69071.871
$ws.Range("M136").Value = -2514.857400000001
$ws.Range("N136").Value = -74171.871

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3086.3157
$ws.Range("I3").Value = 1890.625
$ws.Range("K3").Value = 5671.875
$ws.Range("M3").Value = -5559.875

$ws.Range("I5").Value = 776.6667
$ws.Range("J5").Value = 2045.909
$ws.Range("K5").Value = 2330.0001
$ws.Range("L5").Value = 6137.727000000001
$ws.Range("M5").Value = -2218.0001
$ws.Range("N5").Value = -6361.727000000001

$ws.Range("H68").Value = 11789.444
$ws.Range("I68").Value = 474.75
$ws.Range("J68").Value = 20841.2
$ws.Range("K68").Value = 1424.25
$ws.Range("L68").Value = 62523.60000000001
$ws.Range("M68").Value = -613.25
$ws.Range("N68").Value = -64145.60000000001

$ws.Range("H71").Value = 11789.444
$ws.Range("I71").Value = 474.75
$ws.Range("J71").Value = 20841.2
$ws.Range("K71").Value = 4272.75
$ws.Range("L71").Value = 187570.8
$ws.Range("M71").Value = -216.75
$ws.Range("N71").Value = -195682.8

$ws.Range("H113").Value = 801
$ws.Range("I113").Value = 602.75
$ws.Range("J113").Value = 970.9286
$ws.Range("K113").Value = 1808.25
$ws.Range("L113").Value = 2912.7858
$ws.Range("M113").Value = 361.75
$ws.Range("N113").Value = -7252.7858

$ws.Range("H131").Value = 656.1778
$ws.Range("J131").Value = 696.2436
$ws.Range("L131").Value = 2088.7308
$ws.Range("N131").Value = -12168.7308

$ws.Range("I135").Value = 776.6667
$ws.Range("J135").Value = 2045.909
$ws.Range("K135").Value = 6990.0003
$ws.Range("L135").Value = 18413.181
$ws.Range("M135").Value = -4455.0003
$ws.Range("N135").Value = -23483.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1000
$ws.Range("J6").Value = 1000
$ws.Range("L6").Value = 1000
$ws.Range("N6").Value = -1226

$ws.Range("H16").Value = 1000
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1500

$ws.Range("H132").Value = 3545180.8
$ws.Range("I132").Value = 4382793
$ws.Range("J132").Value = 75073.14
$ws.Range("K132").Value = 13148379
$ws.Range("L132").Value = 225219.42
$ws.Range("M132").Value = -13145849
$ws.Range("N132").Value = -230279.42

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 21992.5
$ws.Range("J58").Value = 21992.5
$ws.Range("L58").Value = 21992.5
$ws.Range("N58").Value = -22608.5

$ws.Range("H132").Value = 19232282
$ws.Range("I132").Value = 41667696
$ws.Range("J132").Value = 1927.2858
$ws.Range("K132").Value = 125003088
$ws.Range("L132").Value = 5781.857400000001
$ws.Range("M132").Value = -125000558
$ws.Range("N132").Value = -10841.8574
